$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.070.09'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '3.264.63'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.95'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.12'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  +1.17%  '
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.62'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('D12').Value = '3.835.84'
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.50'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '68.061.19'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('E16').Value = '  -1.50%  '
$ws.Range('D17').Value = '3.289.19'
$ws.Range('E17').Value = '  +2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.71'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '416.13'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +7.12%  '
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.44'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.42'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.01'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.94'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('E33').Value = '  -3.02%  '
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '162.84'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.44'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.04'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('E39').Value = '  -2.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.46'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.80%  '
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('D42').Value = '2.637.87'
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.79'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.43'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.42%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0675'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '337.11'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.33'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.976'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('E51').Value = '  -1.46%  '
